$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.673.52'
$ws.Range("E2").Value = '  -1.49%  '
$ws.Range("D3").Value = '2.674.62'
$ws.Range("E3").Value = '  -1.25%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '598.15'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -0.26%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '166.28'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +3.08%  '
$ws.Range("E7").Value = '  +0.00%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.547'
$c.Style = "Normal"
$ws.Range("E8").Value = '  +0.47%  '
$ws.Range("D9").Value = '2.674.26'
$ws.Range("E9").Value = '  -1.24%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.145'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +2.22%  '
$ws.Range("E11").Value = '  +1.16%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.361'
$c.Style = "Normal"
$ws.Range("E13").Value = '  -1.60%  '
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '27.87'
$c.Style = "Normal"
$ws.Range("E14").Value = '  -1.95%  '
$ws.Range("D15").Value = '3.162.76'
$ws.Range("E15").Value = '  -0.98%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '0.0000186'
$c.Style = "Normal"
$ws.Range("E16").Value = '  -1.62%  '
$ws.Range("D17").Value = '67.389.97'
$ws.Range("E17").Value = '  -1.72%  '
$ws.Range("D18").Value = '2.683.61'
$ws.Range("E18").Value = '  -0.96%  '
$ws.Range("E19").Value = '  -0.95%  '
$ws.Range("E20").Value = '  +0.73%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '364.00'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -0.32%  '
$ws.Range("E22").Value = '  -3.48%  '
$ws.Range("E23").Value = '  -1.58%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '2.04'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -3.91%  '
$ws.Range("E25").Value = '  -0.05%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '70.99'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -4.28%  '
$ws.Range("E27").Value = '  +1.34%  '
$ws.Range("D28").Value = '2.820.97'
$ws.Range("E28").Value = '  -0.62%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '0.0000103'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -2.76%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -0.07%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '557.52'
$c.Style = "Normal"
$ws.Range("E31").Value = '  -6.48%  '
$ws.Range("E32").Value = '  -3.15%  '
$ws.Range("E33").Value = '  -3.64%  '
$ws.Range("E34").Value = '  -1.08%  '
$ws.Range("E35").Value = '  -1.56%  '
$ws.Range("E36").Value = '  +0.01%  '
$ws.Range("E37").Value = '  -4.76%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '19.55'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -1.64%  '
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '156.00'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -3.38%  '
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.374'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -1.89%  '
$ws.Range("E41").Value = '  -1.81%  '
$ws.Range("E42").Value = '  -4.36%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '17.94'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -0.37%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '2.53'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -5.83%  '
$ws.Range("E45").Value = '  +0.03%  '
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '40.35'
$c.Style = "Normal"
$ws.Range("E46").Value = '  -0.92%  '
$ws.Range("E47").Value = '  -5.51%  '
$ws.Range("E48").Value = '  -2.48%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '153.62'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -2.86%  '
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '3.86'
$c.Style = "Normal"
$ws.Range("E50").Value = '  -2.13%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '1.73'
$c.Style = "Normal"
$ws.Range("E51").Value = '  -2.88%  '
